$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.387.40'
$ws.Range("E2").Value = '  -2.70%  '

$ws.Range("D3").Value = '1.986.09'
$ws.Range("E3").Value = '  -1.29%  '

$ws.Range("E4").Value = '  -0.08%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '235.99'
$c.ClearFormats()
$ws.Range("E5").Value = '  -8.67%  '

$ws.Range("E6").Value = '  -3.26%  '

$ws.Range("E7").Value = '  +0.05%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '54.45'
$c.ClearFormats()
$ws.Range("E8").Value = '  -2.69%  '

$ws.Range("E9").Value = '  -4.28%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '58.07'
$c.ClearFormats()
$ws.Range("E10").Value = '  +2.53%  '

$ws.Range("E11").Value = '  -2.89%  '

$ws.Range("E12").Value = '  -3.01%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '14.16'
$c.ClearFormats()
$ws.Range("E13").Value = '  -0.63%  '

$ws.Range("D14").Value = '2.280.10'
$ws.Range("E14").Value = '  -1.35%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '20.17'
$c.ClearFormats()
$ws.Range("E15").Value = '  -2.77%  '

$ws.Range("E16").Value = '  -6.18%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '5.05'
$c.ClearFormats()
$ws.Range("E17").Value = '  -3.92%  '

$ws.Range("D18").Value = '1.988.15'
$ws.Range("E18").Value = '  -0.87%  '

$ws.Range("D19").Value = '36.362.53'
$ws.Range("E19").Value = '  -2.40%  '

$ws.Range("E20").Value = '  -2.67%  '

$ws.Range("D21").Value = '0.0₃0803'
$ws.Range("E21").Value = '  -4.42%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.27'
$c.ClearFormats()
$ws.Range("E22").Value = '  +2.04%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '221.28'
$c.ClearFormats()
$ws.Range("E23").Value = '  -2.98%  '

$ws.Range("E24").Value = '  -0.08%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.37'
$c.ClearFormats()
$ws.Range("E25").Value = '  +1.06%  '

$ws.Range("E26").Value = '  -9.42%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '162.84'
$c.ClearFormats()
$ws.Range("E27").Value = '  -0.85%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '8.66'
$c.ClearFormats()
$ws.Range("E28").Value = '  -3.57%  '

$ws.Range("E29").Value = '  -0.77%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '18.79'
$c.ClearFormats()
$ws.Range("E30").Value = '  -4.43%  '

$ws.Range("E31").Value = '  +0.80%  '

$ws.Range("E32").Value = '  -2.94%  '

$ws.Range("E33").Value = '  -5.64%  '

$ws.Range("E34").Value = '  -6.61%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.ClearFormats()
$ws.Range("E35").Value = '  -6.49%  '

$ws.Range("E36").Value = '  -2.52%  '

$ws.Range("E37").Value = '  -0.09%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.31'
$c.ClearFormats()
$ws.Range("E38").Value = '  -0.90%  '

$ws.Range("B39").Value = 'WEMIXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.76'
$c.ClearFormats()
$ws.Range("E39").Value = '  -3.62%  '

$ws.Range("E40").Value = '  +5.53%  '

$ws.Range("E41").Value = '  -1.49%  '

$ws.Range("D42").Value = '1.450.54'
$ws.Range("E42").Value = '  +4.41%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0918'
$c.ClearFormats()
$ws.Range("E43").Value = '  -2.17%  '

$ws.Range("E44").Value = '  -5.16%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '89.43'
$c.ClearFormats()
$ws.Range("E45").Value = '  -0.68%  '

$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '1.10'
$c.ClearFormats()
$ws.Range("E46").Value = '  -9.21%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '14.96'
$c.ClearFormats()
$ws.Range("E47").Value = '  -4.43%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.ClearFormats()
$ws.Range("E48").Value = '  -2.79%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.88'
$c.ClearFormats()
$ws.Range("E49").Value = '  -0.70%  '

$ws.Range("E50").Value = '  -3.79%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '3.67'
$c.ClearFormats()
$ws.Range("E51").Value = '  +7.09%  '
